$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 2007年 row (row 2) was removed from the data table; all subsequent
# rows (2010年, 2012年, 2015年, 2017年) shift up by one row.
$ws.Rows(2).Delete()
